# Apply the cryptos-list price/volume refresh described by the commit.
# Column D holds price strings that look numeric ("21.94", "0.0582", ...);
# Excel auto-converts such literals to numbers on assignment, so each is
# written with a leading apostrophe (Excel's text-prefix convention) to keep
# it text, matching the original (non-numeric) inlineStr/shared-string cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.796.70"
$ws.Cells.Item(2, 5).Value = "  -2.45%  "
$ws.Cells.Item(3, 4).Value = "'1.560.95"
$ws.Cells.Item(3, 5).Value = "  -0.48%  "
$ws.Cells.Item(4, 5).Value = "  +0.20%  "
$ws.Cells.Item(5, 4).Value = "'205.75"
$ws.Cells.Item(5, 5).Value = "  -0.99%  "
$ws.Cells.Item(6, 5).Value = "  -2.07%  "
$ws.Cells.Item(7, 5).Value = "  +0.17%  "
$ws.Cells.Item(8, 4).Value = "'21.94"
$ws.Cells.Item(8, 5).Value = "  -0.21%  "
$ws.Cells.Item(9, 5).Value = "  -0.35%  "
$ws.Cells.Item(10, 4).Value = "'0.0582"
$ws.Cells.Item(10, 5).Value = "  -1.22%  "
$ws.Cells.Item(11, 5).Value = "  -0.22%  "
$ws.Cells.Item(12, 4).Value = "'1.786.53"
$ws.Cells.Item(12, 5).Value = "  -0.32%  "
$ws.Cells.Item(13, 4).Value = "'1.559.28"
$ws.Cells.Item(13, 5).Value = "  -0.51%  "
$ws.Cells.Item(14, 4).Value = "'3.74"
$ws.Cells.Item(14, 5).Value = "  -2.11%  "
$ws.Cells.Item(15, 4).Value = "'0.512"
$ws.Cells.Item(15, 5).Value = "  -1.36%  "
$ws.Cells.Item(16, 4).Value = "'61.59"
$ws.Cells.Item(16, 5).Value = "  -2.68%  "
$ws.Cells.Item(17, 4).Value = "'26.823.60"
$ws.Cells.Item(17, 5).Value = "  -2.38%  "
$ws.Cells.Item(18, 4).Value = "'215.04"
$ws.Cells.Item(18, 5).Value = "  +0.47%  "
$ws.Cells.Item(19, 4).Value = "'7.31"
$ws.Cells.Item(19, 5).Value = "  +0.69%  "
$ws.Cells.Item(20, 4).Value = "'0.0₃0678"
$ws.Cells.Item(20, 5).Value = "  -1.68%  "
$ws.Cells.Item(21, 5).Value = "  +0.26%  "
$ws.Cells.Item(22, 5).Value = "  -0.89%  "
$ws.Cells.Item(23, 4).Value = "'9.34"
$ws.Cells.Item(24, 5).Value = "  -1.06%  "
$ws.Cells.Item(25, 4).Value = "'151.81"
$ws.Cells.Item(25, 5).Value = "  -0.87%  "
$ws.Cells.Item(26, 4).Value = "'6.73"
$ws.Cells.Item(26, 5).Value = "  -1.10%  "
$ws.Cells.Item(27, 4).Value = "'14.85"
$ws.Cells.Item(27, 5).Value = "  -1.15%  "
$ws.Cells.Item(28, 5).Value = "  +0.24%  "
$ws.Cells.Item(29, 5).Value = "  -1.49%  "
$ws.Cells.Item(30, 4).Value = "'0.0461"
$ws.Cells.Item(30, 5).Value = "  -1.96%  "
$ws.Cells.Item(31, 5).Value = "  -3.91%  "
$ws.Cells.Item(32, 5).Value = "  -1.64%  "
$ws.Cells.Item(33, 4).Value = "'1.388.05"
$ws.Cells.Item(33, 5).Value = "  +1.87%  "
$ws.Cells.Item(34, 5).Value = "  -1.24%  "
$ws.Cells.Item(35, 4).Value = "'1.55"
$ws.Cells.Item(35, 5).Value = "  +1.00%  "
$ws.Cells.Item(36, 5).Value = "  -0.26%  "
$ws.Cells.Item(37, 4).Value = "'0.932"
$ws.Cells.Item(37, 5).Value = "  -4.48%  "
$ws.Cells.Item(39, 4).Value = "'0.808"
$ws.Cells.Item(40, 5).Value = "  -3.88%  "
$ws.Cells.Item(41, 5).Value = "  +0.24%  "
$ws.Cells.Item(42, 5).Value = "  +2.89%  "
$ws.Cells.Item(43, 4).Value = "'5.42"
$ws.Cells.Item(43, 5).Value = "  +2.68%  "
$ws.Cells.Item(44, 2).Value = "MXToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(44, 4).Value = "'2.18"
$ws.Cells.Item(44, 5).Value = "  +1.40%  "
$ws.Cells.Item(45, 2).Value = "RenderToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(45, 4).Value = "'1.77"
$ws.Cells.Item(45, 5).Value = "  -1.50%  "
$ws.Cells.Item(46, 4).Value = "'63.16"
$ws.Cells.Item(46, 5).Value = "  -1.57%  "
$ws.Cells.Item(47, 4).Value = "'1.698.48"
$ws.Cells.Item(47, 5).Value = "  -0.39%  "
$ws.Cells.Item(48, 4).Value = "'85.64"
$ws.Cells.Item(48, 5).Value = "  +0.14%  "
$ws.Cells.Item(49, 4).Value = "'0.0₇0972"
$ws.Cells.Item(49, 5).Value = "  -1.74%  "
$ws.Cells.Item(50, 4).Value = "'0.0493"
$ws.Cells.Item(50, 5).Value = "  -0.33%  "
$ws.Cells.Item(51, 4).Value = "'0.0943"
$ws.Cells.Item(51, 5).Value = "  -1.07%  "
